# Append one new log row to the end of each worksheet (the device log
# captured a new record on every sheet after the last save).
$wb = $excel.ActiveWorkbook

$rows = @(
    @{ Sheet = "ROW50-FE-LIFTER";  Row = 23; A = 45735.62891842593; B = "0x01,0x90";  C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"; D = "0x01,0x7e"; E = "0xe";  F = 400; G = [double]"5.68631262647114e+23";       GIsText = $false; H = 382; I = 14 },
    @{ Sheet = "ROW50-MID-LIFTER"; Row = 25; A = 45735.60466435185; B = "0x01,0x90 "; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; D = "0x01,0x82"; E = "0x19"; F = 400; G = "568631262647113771663628"; GIsText = $true;  H = 386; I = 25 },
    @{ Sheet = "ROW11-FE-LIFTER";  Row = 23; A = 45735.65168590278; B = "0x01,0x90";  C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; D = "0x01,0x7e"; E = "0x14"; F = 400; G = [double]"5.68631262647114e+23";       GIsText = $false; H = 382; I = 20 },
    @{ Sheet = "ROW11-MID-LIFTER"; Row = 23; A = 45735.79907114583; B = "0x01,0x90";  C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; D = "0x01,0x86"; E = "0x19"; F = 400; G = [double]"5.68631262647114e+23";       GIsText = $false; H = 390; I = 25 }
)

foreach ($entry in $rows) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    $r = $entry.Row

    # Column A: timestamp, formatted the same way as the rest of the column.
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 1).Value = $entry.A

    # Columns B-E: raw hex-byte strings.
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E

    # Column F: plain numeric length.
    $ws.Cells.Item($r, 6).Value = $entry.F

    # Column G: decoded ID - keep as text when the magnitude needs full
    # precision (large integer), otherwise store the plain number.
    if ($entry.GIsText) {
        $ws.Cells.Item($r, 7).NumberFormat = "@"
    }
    $ws.Cells.Item($r, 7).Value = $entry.G

    # Columns H-I: numeric decoded length / checksum.
    $ws.Cells.Item($r, 8).Value = $entry.H
    $ws.Cells.Item($r, 9).Value = $entry.I
}
